$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.673.14"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.119.73"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.32%  "
$ws.Range("D5").Value = "'338.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").Value = "'0.5275"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.4559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "'53.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'0.09119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'24.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "2.131.12"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'6.843"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'8.096"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "'98.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "'0.00001174"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("D18").Value = "'1.015"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "'0.06702"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "'19.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "'1.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "'6.450"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "30.747.84"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'12.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").Value = "'2.381"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "2.364.68"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'22.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'165.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").Value = "'2.554"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "'135.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("D31").Value = "'1.210"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'6.410"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "'1.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'3.951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'10.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").Value = "'5.970"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.61%  "
$ws.Range("D38").Value = "'0.02666"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").Value = "'0.06886"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "'12.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'0.6912"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'1.268"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "'15.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.51%  "
$ws.Range("D45").Value = "'0.6486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'0.00000000369"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.29%  "
$ws.Range("D48").Value = "'3.706"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").Value = "'1.259"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").Value = "'83.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'0.07311"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.59%  "
